# Navigation UI upgrades and edit page UI fix
#
# Updates the answer-choice text for two questions on the MultipleChoice
# sheet ("fruit" and "favourite place") and fixes a typo on the pasta
# question, then moves the saved selection on that sheet.
#
# NOTE: cell values are written in the same left-to-right/right-to-left
# order the original edit used so that any newly introduced strings land
# in the same position in the workbook's shared string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MultipleChoice")

# Row 2 - "Which of the following is a fruit?" answer choices
$ws.Range("E2").Value = "Applesssssss"
$ws.Range("D2").Value = "Basilsssssss"
$ws.Range("C2").Value = "Beefsssssss"
$ws.Range("B2").Value = "Porksssssssssss"
$ws.Range("F2").Value = "None of the Abovesss"

# Row 3 - "Which is your favourite place?" answer choices
$ws.Range("B3").Value = "Parkssssssssss"
$ws.Range("C3").Value = "Homesssssss"
$ws.Range("D3").Value = "Mallssssssss"
$ws.Range("E3").Value = "Poolssssssss"
$ws.Range("F3").Value = "Theatersssssss"

# Row 11 - "Which is a type of pasta?" typo fix
$ws.Range("D11").Value = "Fried Ric"

# Move the active selection on the MultipleChoice tab
$ws.Activate()
$ws.Range("D12").Select()
